# Weekly update: insert the next day's Ciboulette price record as a new
# row 97 (pushing the existing rows 97:189 down to 98:190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 97:189 down one row, leaving a blank (but formatted) row 97.
$ws.Rows("97:97").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = "Femacal de La Calera"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44484
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = 100112039
$ws.Range("G97").Value = "Ciboulette"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 160
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 1500
$ws.Range("N97").Value = "`$/docena de atados"
$ws.Range("O97").Value = "Provincia de Quillota"
$ws.Range("P97").Value = 500
$ws.Range("Q97").Value = 3
$ws.Range("R97").Value = "Hortaliza"
